# "Generate Report for Handoff" — the localization-status CI report was
# regenerated; the row for b9d37f6d-eb4e-42b1-8367-13399e1fe118.md (which is
# "Ready for handoff") picked up a fresh "Latest HO Xliff Generate Date"
# timestamp on the Overview sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G5").Value = "2016-09-05 22:53:12"
